$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage the new (pretty-printed JSON) text on a scratch cell far away so that
# Excel's auto row-height-on-edit (triggered by the many embedded newlines)
# lands on a throwaway row instead of a row we are keeping.
$ws.Range("A100").Value = 'questions = [
    {
        "title": "You are the HR manager for a startup software company. You have recently implemented a $20,000 employment branding campaign to spotlight the company\u2019s culture and benefits to potential job seekers.  Which Key Performance Indicator (KPI) should you use to measure the financial gain of the campaign?",
        "ques_type": 2,
        "options": [
            "Return on investment.",
            "Turnover rate.",
            "Position time to fill.",
            "New hire retention rate."
        ],
        "score": "Return on investment."
    },
    {
        "title": "You are the HR manager for an auto manufacturer. The accounting manager has confided to you that she is nervous about leading an upcoming termination meeting because she is unsure how to deal with the significant pushback that is expected from the separating employee. How should you offer to help her?",
        "ques_type": 2,
        "options": [
            "\u201cWould you like me to conduct the meeting and update you afterward?\u201d",
            "\u201cWould you like me to cancel the meeting and give the employee another opportunity to improve?\u201d",
            "\u201cWould you like to role-play the conversation with me in advance and look at some talking points you can use?\u201d",
            "\u201cWould you like me to authorize an increase in their severance package you can use if they push back?\u201d"
        ],
        "score": "\u201cWould you like to role-play the conversation with me in advance and look at some talking points you can use?\u201d"
    },
    {
        "title": "You are working as the HR manager at a paper company. A warehouse employee who started today has just told you that she is four months pregnant. She said she can only lift 10 pounds (approx. 4.5 kg) instead of the required 50 pounds (approx. 22.5 kg) and has requested a job accommodation. What action should you take?",
        "ques_type": 2,
        "options": [
            "Ask the employee why she did not disclose her pregnancy during the interview.",
            "Terminate the employee without addressing the accommodation.",
            "Change the job description to reduce the lifting requirement.",
            "Determine if there is a position available within the company that can accommodate the lifting restriction."
        ],
        "score": "Determine if there is a position available within the company that can accommodate the lifting restriction."
    },
    {
        "title": "Your chief executive officer has requested that the accounting department consider hiring his nephew as an assistant accountant.  True or false: You should warn him that this may result in a conflict of interest.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "True"
    }
]'

# Move the staged text onto A2 via cut/paste -- this does not re-trigger the
# auto-height calculation, so row 2 keeps its default (non-custom) height.
$ws.Range("A100").Cut($ws.Range("A2"))

# Drop the now-empty scratch row entirely so no stray formatting remains.
$ws.Rows.Item(100).Delete()

# Finally remove the helper row 1 (value 0 / bold+border style); this shifts
# the updated A2 up to A1, matching the target layout.
$ws.Rows.Item(1).Delete()
